$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.269.25"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "2.764.56"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'583.39"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "'160.18"
$ws.Range("E6").Value = "  +4.72%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").Value = "'5.92"
$ws.Range("E10").Value = "  -11.28%  "
$ws.Range("D11").Value = "'0.392"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "3.253.62"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").Value = "'27.03"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").Value = "64.127.29"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "2.769.65"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "'12.29"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("D20").Value = "'363.60"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").Value = "'0.573"
$ws.Range("E22").Value = "  +7.11%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").Value = "'66.71"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("E25").Value = "  +3.79%  "
$ws.Range("D26").Value = "'8.67"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").Value = "0.0₃0945"
$ws.Range("E28").Value = "  +5.77%  "
$ws.Range("D29").Value = "'2.02"
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").Value = "'7.17"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "'1.27"
$ws.Range("E31").Value = "  +5.17%  "
$ws.Range("D32").Value = "'170.49"
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "'20.65"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "'5.02"
$ws.Range("E35").Value = "  +4.26%  "
$ws.Range("E36").Value = "  +2.62%  "
$ws.Range("D37").Value = "'1.84"
$ws.Range("E37").Value = "  +3.54%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "'4.25"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "'336.84"
$ws.Range("E40").Value = "  -2.98%  "
$ws.Range("D41").Value = "'6.19"
$ws.Range("E41").Value = "  +10.97%  "
$ws.Range("D42").Value = "'39.71"
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("D43").Value = "'22.25"
$ws.Range("E43").Value = "  +1.75%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'22.26"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0602"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.646"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0260"
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").Value = "'137.07"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("D49").Value = "'0.103"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'11.06"
$ws.Range("E51").Value = "  +0.92%  "
